# Removal of duplicated test case from todo test case documentation
#
# Row 7 ("Test successful display of to-do items based on user") was a
# duplicate test case. It is removed by shifting up the test-case details
# (columns B:G) of the rows below it, and clearing out the now-trailing row.
# The numbering in column A (2.5, 2.6, 2.7 ...) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the content of row 8 into row 7, and row 9 into row 8.
$ws.Range("B8:G8").Copy() | Out-Null
$ws.Range("B7:G7").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Range("B9:G9").Copy() | Out-Null
$ws.Range("B8:G8").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Application.CutCopyMode = $false

# Clear out the now-duplicated trailing row entirely.
$ws.Range("A9:G9").Clear()

# Leave the active cell where the user ended up after clearing the row.
$ws.Range("A9").Select() | Out-Null
